# Update countries & provincias Spain
# Applies the 21-Jun-2020 07:31 data refresh to the "Pais" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp note -------------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 21 de Junio de 2020 a las 07:31"

# --- Plain numeric refreshes (row order unchanged) ------------------------
# Kazajistan
$ws.Range("B56").Value = 17225
$ws.Range("C56").Value = 446
$ws.Range("D56").Value = 10671
$ws.Range("E56").Value = 6436

# Uzbekistan
$ws.Range("B76").Value = 6216
$ws.Range("C76").Value = 63
$ws.Range("D76").Value = 4290
$ws.Range("E76").Value = 1907

# Surinam
$ws.Range("B160").Value = 302
$ws.Range("E160").Value = 220

# --- Rows whose country label swaps with its neighbour because the new
#     totals changed the ranking order --------------------------------

# Rows 80/81 : Republica de Macedonia <-> Haiti
$ws.Range("A80").Value = "Haiti"
$ws.Range("B80").Value = 5077
$ws.Range("C80").Value = 97
$ws.Range("D80").Value = 24
$ws.Range("E80").Value = 4965
$ws.Range("F80").Value = 0
$ws.Range("G80").Value = 1
$ws.Range("H80").Value = 88

$ws.Range("A81").Value = "Republica de Macedonia"
$ws.Range("B81").Value = 5005
$ws.Range("C81").Value = 0
$ws.Range("D81").Value = 1904
$ws.Range("E81").Value = 2868
$ws.Range("F81").Value = 0
$ws.Range("G81").Value = 0
$ws.Range("H81").Value = 233

# Rows 94/95 : Tailandia <-> Kirguistan
$ws.Range("A94").Value = "Kirguistan"
$ws.Range("B94").Value = 3151
$ws.Range("C94").Value = 170
$ws.Range("D94").Value = 2011
$ws.Range("E94").Value = 1103
$ws.Range("F94").Value = 0
$ws.Range("G94").Value = 2
$ws.Range("H94").Value = 37

$ws.Range("A95").Value = "Tailandia"
$ws.Range("B95").Value = 3148
$ws.Range("C95").Value = 1
$ws.Range("D95").Value = 3018
$ws.Range("E95").Value = 72
$ws.Range("F95").Value = 0
$ws.Range("G95").Value = 0
$ws.Range("H95").Value = 58

# Rows 202/203 : Dominica <-> Fiyi (figures identical, only label swaps)
$ws.Range("A202").Value = "Fiyi"
$ws.Range("A203").Value = "Dominica"

# Rows 208/209 : Islas Turcas y Caicos <-> Santa Sede
$ws.Range("A208").Value = "Santa Sede"
$ws.Range("D208").Value = 12
$ws.Range("H208").Value = 0

$ws.Range("A209").Value = "Islas Turcas y Caicos"
$ws.Range("D209").Value = 11
$ws.Range("H209").Value = 1

# Rows 213/214 : Islas Virgenes Britanicas <-> Papua Nueva Guinea
$ws.Range("A213").Value = "Papua Nueva Guinea"
$ws.Range("D213").Value = 8
$ws.Range("H213").Value = 0

$ws.Range("A214").Value = "Islas Virgenes Britanicas"
$ws.Range("D214").Value = 7
$ws.Range("H214").Value = 1
